$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (12 data rows total, written into rows 2..13), header row 1
# stays blank/unchanged.
$data = @(
    @("10004394", "MENTOS FRUIT ROLL 37", "TH1MKT", "1", "1", "RT,(E-1B)"),
    @("10030201", "MENTOS CND ROL AGR37", "TH1MKT", "1", "2", "RT,(E-1B)"),
    @("20137700", "MAYASI CSHW KTCY 35G", "TH1MKT", "2", "1", "RT,(E-1B)"),
    @("20103453", "SAORI SAUS LD HTM133", "TH1MKT", "2", "2", "RT,(E-1B)"),
    @("20009966", "SAORI SAUS TIRAM 133", "TH1MKT", "2", "3", "RT,(E-1B)"),
    @("20009973", "SAORI SAUS TERIYK135", "TH1MKT", "2", "4", "RT,(E-1B)"),
    @("20139596", "PPSODNT SENSI EXP 60", "TH1MKT", "3", "1", "PT,(E-1B)"),
    @("10040202", "SOFFELL A.NYMK K/J60", "TH1MKT", "3", "3", "RT,(E-3.5B)"),
    @("20134253", "SO FRSH HOT 2X10ML", "TH1MKT", "3", "4", "RT,(E-1B)"),
    @("20040313", "MY BABY TELON PLUS30", "TH1MKT", "3", "5", "RT,(E-3B)"),
    @("20122879", "RNSO MLTO KRN STR510", "TH1MKT", "4", "1", "PT"),
    @("20140001", "RINSO PURE LIQ 510G", "TH1MKT", "4", "2", "PT")
)

# Sheet grows from 12 data rows (A1:F12) to 13 data rows (A1:F13): copy the
# last existing data row's formatting (border style) down into new row 13
# before writing into it, so it matches the look of the other rows.
$ws.Range("A12:F12").Copy() | Out-Null
$ws.Range("A13:F13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Columns A, D and E hold digit-only values (item codes / numbering) that
# must stay TEXT (as in the source file) rather than being auto-coerced to
# numbers, so mark those ranges as Text before writing the values.
$ws.Range("A2:A13").NumberFormat = "@"
$ws.Range("D2:E13").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowVals[$c]
    }
}

# Column F widened slightly (11 -> 13) per the updated layout.
# (12.16 is the ColumnWidth input that round-trips to a saved width of 13
# in the generated XML, matching the target exactly.)
$ws.Columns.Item(6).ColumnWidth = 12.16
